# "Fixed typos in 'important social science issues' slide" + table style update
# (Add files via upload)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 15: update the Machine Learning vs Social Science table's style
# ---------------------------------------------------------------------------
$sTable = $p.Slides.Item(15)
$tbl = $sTable.Shapes.Item(2).Table
$tbl.ApplyStyle("{9D89CBAC-87FF-44C3-AA51-7410B9017C42}")

# ---------------------------------------------------------------------------
# 2) Slide 17: "Important Social Science Issues"
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(17)

$titleShape = $s.Shapes.Item(1)
$bodyShape = $s.Shapes.Item(2)

# Nudge the title & body placeholders up slightly (EMU -129650 / 356425,
# expressed in points - the COM layer stores Top/Left as single-precision
# floats, so these literals are chosen to land on the exact target EMU).
$titleShape.Top = -10.208662033081055
$bodyShape.Top = 28.06496238708496

$tr = $bodyShape.TextFrame.TextRange

# "Too many rows" bullet sub-point was mislabeled; reword it.
$tr.Paragraphs(7).Text = "p hacking, in-sample fitting and testing"

# Add the missing "Post-treatment controls" sub-bullet under "Too many columns".
$tooManyColumns = $tr.Paragraphs(10)
$noise = $tr.Paragraphs(11)
$noise.InsertBefore("Post-treatment controls`r")
$newBullet = $tr.Paragraphs(11)
$newBullet.IndentLevel = $tooManyColumns.IndentLevel + 1
